$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# (Shared by Overview!B2/C2/B3/C3 and the Status column on each language sheet.)
$statusText = "Handed back: in sync with en-US"
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# --- Latest Handback DateTime ---
$wsZh.Range("H2").Value = "2016-03-13 00:55:26"
$wsZh.Range("H3").Value = "2016-03-13 00:55:26"
$wsDe.Range("H2").Value = "2016-03-13 00:55:32"
$wsDe.Range("H3").Value = "2016-03-13 00:55:32"

# --- New "Latest Target File" / "Latest Handback File" columns (F, G) ---
# zh-cn sheet, rows 2 and 3
$wsZh.Range("F2").Value = "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md"
$wsZh.Range("G2").Value = "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.zh-cn.xlf"
$wsZh.Range("F3").Value = "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md"
$wsZh.Range("G3").Value = "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/71ef84107f5b37e069618a30ac588778f1039301/e2e/af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md", "", "", "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/18965593b431071e8578155fc915a7ee46c0f7b2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.zh-cn.xlf", "", "", "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/71ef84107f5b37e069618a30ac588778f1039301/e2e/af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md", "", "", "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/18965593b431071e8578155fc915a7ee46c0f7b2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.zh-cn.xlf", "", "", "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.zh-cn.xlf")

$wsZh.Range("F2").Font.Underline = 2
$wsZh.Range("F2").Font.Color = 15570276
$wsZh.Range("G2").Font.Underline = 2
$wsZh.Range("G2").Font.Color = 15570276
$wsZh.Range("F3").Font.Underline = 2
$wsZh.Range("F3").Font.Color = 15570276
$wsZh.Range("G3").Font.Underline = 2
$wsZh.Range("G3").Font.Color = 15570276

# de-de sheet, rows 2 and 3
$wsDe.Range("F2").Value = "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md"
$wsDe.Range("G2").Value = "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.de-de.xlf"
$wsDe.Range("F3").Value = "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md"
$wsDe.Range("G3").Value = "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/71ef84107f5b37e069618a30ac588778f1039301/e2e/af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md", "", "", "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ac7bc18876f4da0f542fbf6565d41439f04214f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.de-de.xlf", "", "", "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/71ef84107f5b37e069618a30ac588778f1039301/e2e/af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md", "", "", "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ac7bc18876f4da0f542fbf6565d41439f04214f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.de-de.xlf", "", "", "af01bf42-34aa-4cbf-876c-d360ca1fd7b9.104ac376faed4118def4e901bf122f240ac8750d.de-de.xlf")

$wsDe.Range("F2").Font.Underline = 2
$wsDe.Range("F2").Font.Color = 15570276
$wsDe.Range("G2").Font.Underline = 2
$wsDe.Range("G2").Font.Color = 15570276
$wsDe.Range("F3").Font.Underline = 2
$wsDe.Range("F3").Font.Color = 15570276
$wsDe.Range("G3").Font.Underline = 2
$wsDe.Range("G3").Font.Color = 15570276

Write-Host "Generate Report for Handback: done"
